# RoboRIO Ports and other controls.xlsx - apply commit "Updated Ports and other controls file"
#
# Summary of the edit:
#  - Sheet "RoboRIO Ports": a new column D ("Inverts") is introduced, the old
#    "Left and Right Motors" label in C2 is replaced by "RobotDrive", and the
#    encoder / Sonar rows (10, 11, 15) get the built-in "Bad" (red) cell style
#    applied across columns A:C.
#  - Sheet "Controller Map" content is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RoboRIO Ports")

# ---------------------------------------------------------------------------
# Column A (Port) + Column B (Motor/Sensor) + Column C (RefNum (Programming))
# + new Column D (Inverts)
# ---------------------------------------------------------------------------

$table = New-Object 'object[,]' 19,4

# Row 1 - header
$table[0,0] = "Port"
$table[0,1] = "Motor/Sensor"
$table[0,2] = "RefNum (Programming)"
$table[0,3] = "Inverts"

# Row 2
$table[1,0] = "PWM0"
$table[1,1] = "Left Drive Motors"
$table[1,2] = "RobotDrive"
$table[1,3] = "F"

# Row 3
$table[2,0] = "PWM1"
$table[2,1] = "Right Drive Motors"
$table[2,2] = $null
$table[2,3] = "T"

# Row 4
$table[3,0] = "PWM2"
$table[3,1] = "Intake (rubber bands)"
$table[3,2] = "Intake"
$table[3,3] = "F"

# Row 5
$table[4,0] = "PWM3"
$table[4,1] = "Scaler (Window motors)"
$table[4,2] = "Scaler"
$table[4,3] = "F"

# Row 6
$table[5,0] = "PWM4"
$table[5,1] = "Scaler Winch #1"
$table[5,2] = "Scaler Left Lift"
$table[5,3] = "F"

# Row 7
$table[6,0] = "PWM5"
$table[6,1] = "Scaler Winch #2"
$table[6,2] = "Scaler Right Lift"
$table[6,3] = "F"

# Row 8
$table[7,0] = "CAN0"
$table[7,1] = "Intake Arm #1"
$table[7,2] = "Arm 1"
$table[7,3] = "T"

# Row 9
$table[8,0] = "CAN1"
$table[8,1] = "Intake Arm #2"
$table[8,2] = "Arm 2"
$table[8,3] = "F"

# Row 10
$table[9,0] = "DIO 0/1"
$table[9,1] = "Left Wheel Encoder"
$table[9,2] = "Left Encoder"
$table[9,3] = $null

# Row 11
$table[10,0] = "DIO 2/3"
$table[10,1] = "Right Wheel Encoder"
$table[10,2] = "Right Encoder"
$table[10,3] = $null

# Row 12
$table[11,0] = "DIO 4/9"
$table[11,1] = "Beam Break Sensor"
$table[11,2] = "Beam Break"
$table[11,3] = $null

# Row 13
$table[12,0] = "DIO6"
$table[12,1] = "Ball intake sensor"
$table[12,2] = "Ball Switch"
$table[12,3] = $null

# Row 14
$table[13,0] = "AI0"
$table[13,1] = "Potentiometer"
$table[13,2] = "Pot"
$table[13,3] = "T"

# Row 15
$table[14,0] = "AI1"
$table[14,1] = "Sonar"
$table[14,2] = "Sonar"
$table[14,3] = $null

# Row 16
$table[15,0] = "I2C MXP"
$table[15,1] = "NavX"
$table[15,2] = "NavX"
$table[15,3] = $null

# Row 17
$table[16,0] = "USB0"
$table[16,1] = "Left Joystick"
$table[16,2] = "Joystick 0"
$table[16,3] = $null

# Row 18
$table[17,0] = "USB1"
$table[17,1] = "Right Joystick"
$table[17,2] = "Joystick 1"
$table[17,3] = $null

# Row 19
$table[18,0] = "USB2"
$table[18,1] = "Controller"
$table[18,2] = "Operator"
$table[18,3] = $null

$ws.Range("A1:D19").Value = $table

# Column C3 is merged with C2 and must stay blank (it is part of the
# existing C2:C3 merge), so remove the spill-over value that Value= puts there.
$ws.Range("C3").Value = $null

# ---------------------------------------------------------------------------
# Highlight the three "mismatched label" rows (Left/Right Wheel Encoder and
# Sonar) using the built-in "Bad" style, same as Excel's conditional
# formatting gallery would produce.
# ---------------------------------------------------------------------------

$ws.Range("A10:C10").Style = "Bad"
$ws.Range("A11:C11").Style = "Bad"
$ws.Range("A15:C15").Style = "Bad"

# ---------------------------------------------------------------------------
# Column widths / selection to match the refreshed layout
# ---------------------------------------------------------------------------

$ws.Columns.Item(3).ColumnWidth = 23
$ws.Range("D3").Select() | Out-Null
